$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 3.25
$ws.Range("I5").Value = 2.7
$ws.Range("J5").Value = 1.17
$ws.Range("K5").Value = 5
$ws.Range("U5").Value = 13
$ws.Range("W5").Value = 34
$ws.Range("Z5").Value = 5
$ws.Range("AF5").Value = 11
$ws.Range("AG5").Value = 12
$ws.Range("AI5").Value = 29

# Row 6
$ws.Range("G6").Value = 3.5
$ws.Range("I6").Value = 2.3
$ws.Range("T6").Value = 8.5
$ws.Range("W6").Value = 41
$ws.Range("AE6").Value = 6
$ws.Range("AH6").Value = 21

# Row 7
$ws.Range("G7").Value = 1.62
$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 1.03
$ws.Range("K7").Value = 15
$ws.Range("L7").Value = 1.2
$ws.Range("M7").Value = 4.33
$ws.Range("W7").Value = 12
$ws.Range("AI7").Value = 41

# Row 8
$ws.Range("G8").Value = 1.38
$ws.Range("J8").Value = 1.03
$ws.Range("K8").Value = 17
$ws.Range("N8").Value = 1.5
$ws.Range("O8").Value = 2.5
$ws.Range("P8").Value = 1.25
$ws.Range("Q8").Value = 3.75
$ws.Range("R8").Value = 1.73
$ws.Range("S8").Value = 2
$ws.Range("T8").Value = 9.5
$ws.Range("U8").Value = 8
$ws.Range("V8").Value = 8.5
$ws.Range("Z8").Value = 17
$ws.Range("AA8").Value = 9
$ws.Range("AE8").Value = 23

# Row 9
$ws.Range("G9").Value = 1.95

# Row 10
$ws.Range("G10").Value = 2.25
$ws.Range("I10").Value = 2.8
$ws.Range("W10").Value = 23
$ws.Range("AE10").Value = 12

# Row 13
$ws.Range("G13").Value = 2.15

# Row 14
$ws.Range("G14").Value = 3.05
$ws.Range("H14").Value = 3.05
$ws.Range("I14").Value = 2.2
$ws.Range("N14").Value = 2.02
$ws.Range("O14").Value = 1.62
$ws.Range("P14").Value = 1.4
$ws.Range("Q14").Value = 2.42
$ws.Range("T14").Value = 7.5
$ws.Range("U14").Value = 13
$ws.Range("V14").Value = 9.25
$ws.Range("W14").Value = 32
$ws.Range("X14").Value = 22
$ws.Range("Y14").Value = 29
$ws.Range("Z14").Value = 8
$ws.Range("AA14").Value = 5.2
$ws.Range("AB14").Value = 11.75
$ws.Range("AC14").Value = 50
$ws.Range("AD14").Value = 400
$ws.Range("AE14").Value = 6
$ws.Range("AF14").Value = 8.5
$ws.Range("AG14").Value = 7.6
$ws.Range("AH14").Value = 17.5
$ws.Range("AI14").Value = 15.5
$ws.Range("AJ14").Value = 24

# Row 17
$ws.Range("G17").Value = 2.15
$ws.Range("H17").Value = 3.1
$ws.Range("I17").Value = 3.75
$ws.Range("J17").Value = 1.07
$ws.Range("K17").Value = 9
$ws.Range("L17").Value = 1.36
$ws.Range("M17").Value = 3
$ws.Range("N17").Value = 2.2
$ws.Range("O17").Value = 1.65
$ws.Range("P17").Value = 1.44
$ws.Range("Q17").Value = 2.63
$ws.Range("R17").Value = 1.95
$ws.Range("S17").Value = 1.8
$ws.Range("V17").Value = 9
$ws.Range("W17").Value = 19
$ws.Range("X17").Value = 19
$ws.Range("Z17").Value = 7.5
$ws.Range("AB17").Value = 15
$ws.Range("AD17").Value = 351
$ws.Range("AE17").Value = 9.5

# Row 20
$ws.Range("G20").Value = 3.25
$ws.Range("H20").Value = 3.2
$ws.Range("I20").Value = 2.12
$ws.Range("L20").Value = 1.4
$ws.Range("M20").Value = 2.52
$ws.Range("N20").Value = 2.15
$ws.Range("O20").Value = 1.55
$ws.Range("P20").Value = 1.47
$ws.Range("Q20").Value = 2.32
$ws.Range("R20").Value = 1.93
$ws.Range("S20").Value = 1.7
$ws.Range("T20").Value = 8.25
$ws.Range("U20").Value = 15.5
$ws.Range("V20").Value = 12
$ws.Range("W20").Value = 45
$ws.Range("X20").Value = 32
$ws.Range("Y20").Value = 50
$ws.Range("Z20").Value = 7.7
$ws.Range("AA20").Value = 6.2
$ws.Range("AB20").Value = 17.5
$ws.Range("AC20").Value = 100
$ws.Range("AE20").Value = 6.3
$ws.Range("AF20").Value = 9.25
$ws.Range("AG20").Value = 9.25
$ws.Range("AH20").Value = 19.5
$ws.Range("AI20").Value = 19.5
$ws.Range("AJ20").Value = 37

# Row 21
$ws.Range("G21").Value = 2.32
$ws.Range("H21").Value = 3
$ws.Range("I21").Value = 3.1
$ws.Range("L21").Value = 1.4
$ws.Range("M21").Value = 2.52
$ws.Range("N21").Value = 2.15
$ws.Range("O21").Value = 1.55
$ws.Range("P21").Value = 1.5
$ws.Range("Q21").Value = 2.25
$ws.Range("R21").Value = 1.85
$ws.Range("S21").Value = 1.75
$ws.Range("T21").Value = 6.7
$ws.Range("U21").Value = 10.5
$ws.Range("V21").Value = 9.25
$ws.Range("W21").Value = 24
$ws.Range("X21").Value = 21
$ws.Range("Y21").Value = 35
$ws.Range("Z21").Value = 7.4
$ws.Range("AA21").Value = 5.8
$ws.Range("AB21").Value = 15.5
$ws.Range("AC21").Value = 90
$ws.Range("AD21").Value = 800
$ws.Range("AE21").Value = 7.9
$ws.Range("AF21").Value = 15
$ws.Range("AG21").Value = 11.25
$ws.Range("AH21").Value = 40
$ws.Range("AI21").Value = 30
$ws.Range("AJ21").Value = 45

# Row 23
$ws.Range("G23").Value = 3.1
$ws.Range("I23").Value = 2.15
$ws.Range("J23").Value = 1.05
$ws.Range("K23").Value = 11
$ws.Range("L23").Value = 1.3
$ws.Range("M23").Value = 3.4
$ws.Range("N23").Value = 2
$ws.Range("O23").Value = 1.8
$ws.Range("U23").Value = 15
$ws.Range("V23").Value = 11
$ws.Range("AG23").Value = 9.5
$ws.Range("AH23").Value = 21

# Row 27
$ws.Range("G27").Value = 8.5
$ws.Range("H27").Value = 7
$ws.Range("I27").Value = 1.22
$ws.Range("J27").Value = 23
$ws.Range("K27").Value = 1.02
$ws.Range("L27").Value = 1.1
$ws.Range("M27").Value = 6.5
$ws.Range("N27").Value = 1.33
$ws.Range("O27").Value = 3.25
$ws.Range("P27").Value = 1.2
$ws.Range("Q27").Value = 4.33
$ws.Range("R27").Value = 1.83
$ws.Range("S27").Value = 1.83
$ws.Range("T27").Value = 29
$ws.Range("U27").Value = 51
$ws.Range("V27").Value = 26
$ws.Range("W27").Value = 101
$ws.Range("Z27").Value = 23
$ws.Range("AA27").Value = 15
$ws.Range("AB27").Value = 23
$ws.Range("AD27").Value = 151
$ws.Range("AE27").Value = 11
$ws.Range("AG27").Value = 10
$ws.Range("AH27").Value = 8.5
$ws.Range("AJ27").Value = 23

# Row 31
$ws.Range("K31").Value = 13
$ws.Range("N31").Value = 1.73
$ws.Range("O31").Value = 2.08

# Row 33
$ws.Range("N33").Value = 1.9
$ws.Range("O33").Value = 1.95

# Row 36
$ws.Range("G36").Value = 2.63
$ws.Range("H36").Value = 3
$ws.Range("I36").Value = 2.75
$ws.Range("J36").Value = 1.08
$ws.Range("K36").Value = 8
$ws.Range("U36").Value = 12
$ws.Range("Z36").Value = 8
$ws.Range("AE36").Value = 8.5
$ws.Range("AF36").Value = 13
$ws.Range("AH36").Value = 29

# Row 37
$ws.Range("H37").Value = 3.6
$ws.Range("I37").Value = 4
